$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds dates that must stay as literal text (not be auto-converted
# to Excel date serials). Force a Text number format on the destination
# range before writing the values so the COM input parser keeps them as
# strings.
$ws.Range("C2:C7").NumberFormat = "@"

# Row 2 (overwrite existing expense row)
$ws.Range("A2").Value = 360
$ws.Range("B2").Value = "hair cut"
$ws.Range("C2").Value = "20/12/2026"

# Row 3 (new)
$ws.Range("A3").Value = 390
$ws.Range("B3").Value = "cloth"
$ws.Range("C3").Value = "02/11/2026"

# Row 4 (new)
$ws.Range("A4").Value = 500
$ws.Range("B4").Value = "food"
$ws.Range("C4").Value = "01/11/2026"

# Row 5 (new)
$ws.Range("A5").Value = 7500
$ws.Range("B5").Value = "shopping"
$ws.Range("C5").Value = "28/01/2026"

# Row 6 (new)
$ws.Range("A6").Value = 360
$ws.Range("B6").Value = "cloth"
$ws.Range("C6").Value = "20/01/2026"

# Row 7 (new)
$ws.Range("A7").Value = 50
$ws.Range("B7").Value = "food"
$ws.Range("C7").Value = "01/01/2026"
